$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates: volume/issue number and report date range ---
$ws.Range("A8").Value = "Volume 30   Number  10"
$ws.Range("C9").Value = "Report Covering the Week  3/6/2023  Through  3/12/2023"

# --- Crime statistics table updates (rows 14-29) ---
# For cells that change numeric <-> text type, we first force the new value with a
# leading apostrophe when needed (to stop Excel from auto-parsing numeric-looking text),
# then copy number-format/style from a donor cell in the same column that already has
# the desired resulting style, applying PasteSpecial(xlPasteFormats) so only the format
# (not value) is transferred.

$ws.Range("F14").Value = 2
$ws.Range("G14").Value = "'0"
$ws.Range("D14").Copy() | Out-Null
$ws.Range("G14").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Value = "'***.*"
$ws.Range("E14").Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null
$ws.Range("I14").Value = 2
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = -50
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = -89.473684210526
$ws.Range("I15").Value = 11
$ws.Range("K15").Value = 175
$ws.Range("L15").Value = 83.333333333333
$ws.Range("M15").Value = 266.666666666667
$ws.Range("N15").Value = -8.333333333333
$ws.Range("C16").Value = 15
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 31
$ws.Range("G16").Value = 39
$ws.Range("H16").Value = -20.51282051282
$ws.Range("I16").Value = 83
$ws.Range("J16").Value = 88
$ws.Range("K16").Value = -5.681818181818
$ws.Range("L16").Value = 40.677966101694
$ws.Range("M16").Value = 13.698630136986
$ws.Range("N16").Value = -76.353276353276
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = 15.384615384615
$ws.Range("F17").Value = 56
$ws.Range("G17").Value = 56
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 145
$ws.Range("J17").Value = 126
$ws.Range("K17").Value = 15.079365079365
$ws.Range("L17").Value = 20.833333333333
$ws.Range("M17").Value = 95.945945945946
$ws.Range("N17").Value = -5.228758169934
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = -11.111111111111
$ws.Range("F18").Value = 25
$ws.Range("G18").Value = 33
$ws.Range("H18").Value = -24.242424242424
$ws.Range("I18").Value = 56
$ws.Range("J18").Value = 76
$ws.Range("K18").Value = -26.315789473684
$ws.Range("L18").Value = 86.666666666666
$ws.Range("M18").Value = 40
$ws.Range("N18").Value = -79.259259259259
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = -41.176470588235
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = -13.043478260869
$ws.Range("I19").Value = 125
$ws.Range("J19").Value = 121
$ws.Range("K19").Value = 3.305785123966
$ws.Range("L19").Value = -4.580152671755
$ws.Range("M19").Value = 83.823529411764
$ws.Range("N19").Value = 21.35922330097
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 26
$ws.Range("G20").Value = 36
$ws.Range("H20").Value = -27.777777777777
$ws.Range("I20").Value = 63
$ws.Range("J20").Value = 77
$ws.Range("K20").Value = -18.181818181818
$ws.Range("L20").Value = 125
$ws.Range("M20").Value = 186.363636363636
$ws.Range("N20").Value = -49.6
$ws.Range("C21").Value = 54
$ws.Range("D21").Value = 57
$ws.Range("E21").Value = -5.263157894736
$ws.Range("F21").Value = 183
$ws.Range("G21").Value = 211
$ws.Range("H21").Value = -13.270142180094
$ws.Range("I21").Value = 485
$ws.Range("J21").Value = 494
$ws.Range("K21").Value = -1.821862348178
$ws.Range("L21").Value = 28.306878306878
$ws.Range("M21").Value = 72.597864768683
$ws.Range("N21").Value = -53.049370764762
$ws.Range("D22").Value = 3
$ws.Range("C16").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = -100
$ws.Range("H16").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = -77.777777777777
$ws.Range("J22").Value = 25
$ws.Range("K22").Value = -72
$ws.Range("C23").Value = 9
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = 50
$ws.Range("G23").Value = 23
$ws.Range("H23").Value = 26.086956521739
$ws.Range("I23").Value = 80
$ws.Range("J23").Value = 59
$ws.Range("K23").Value = 35.593220338983
$ws.Range("L23").Value = 31.147540983606
$ws.Range("M23").Value = 60
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -30
$ws.Range("F24").Value = 117
$ws.Range("H24").Value = -13.970588235294
$ws.Range("I24").Value = 301
$ws.Range("J24").Value = 270
$ws.Range("K24").Value = 11.481481481481
$ws.Range("L24").Value = 8.273381294964
$ws.Range("M24").Value = 39.351851851851
$ws.Range("C25").Value = 25
$ws.Range("D25").Value = 21
$ws.Range("E25").Value = 19.047619047619
$ws.Range("F25").Value = 91
$ws.Range("G25").Value = 95
$ws.Range("H25").Value = -4.210526315789
$ws.Range("I25").Value = 197
$ws.Range("J25").Value = 180
$ws.Range("K25").Value = 9.444444444444
$ws.Range("L25").Value = 25.477707006369
$ws.Range("M25").Value = 7.065217391304
$ws.Range("C26").Value = 1
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 6
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 20
$ws.Range("I26").Value = 15
$ws.Range("J26").Value = 8
$ws.Range("K26").Value = 87.5
$ws.Range("L26").Value = 50
$ws.Range("C27").Value = 3
$ws.Range("E27").Value = 50
$ws.Range("I27").Value = 22
$ws.Range("J27").Value = 15
$ws.Range("K27").Value = 46.666666666666
$ws.Range("L27").Value = 29.411764705882
$ws.Range("D28").Value = "'0"
$ws.Range("C28").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value = "'***.*"
$ws.Range("C28").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("G28").Value = 1
$ws.Range("N28").Value = -97.5
$ws.Range("D29").Value = "'0"
$ws.Range("C28").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("E29").Value = "'***.*"
$ws.Range("C28").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null
$ws.Range("G29").Value = 1
$ws.Range("N29").Value = -96.969696969697
